{"js": "// The phrase \" includes the phases as mentioned below-\" (spanning three\n// separate runs) is replaced with \" includes the phases -\", effectively\n// removing the previously added \" as mentioned below\" text and merging\n// the surrounding text back into a single run.\nconst body = context.document.body;\n\nconst results = body.search(\" includes the phases as mentioned below-\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\" includes the phases -\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The phrase \" as mentioned below\" (its own run) between \" includes the\n# phases \" and \"-\" is removed, so the text reads \"... includes the\n# phases -\" and the surrounding runs merge back together.\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$find = $range.Find\n$find.Text = \"as mentioned below\"\n$found = $find.Execute()\n\nif ($found) {\n    $range.Delete()\n}\n"}
